$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new daily status row (row 28) below the existing data (last row was 27)
$ws.Range("A28").Value = 45975
$ws.Range("A28").NumberFormat = "d-mmm-yy"
$ws.Range("B28").Value = 5610
$ws.Range("C28").Value = 4110
$ws.Range("D28").Value = 3790
$ws.Range("E28").Value = 251
$ws.Range("F28").Value = 42
$ws.Range("G28").Value = 25
$ws.Range("H28").Value = 2
$ws.Range("I28").Value = 0

# Match the selection state recorded in the saved workbook (whole new row selected)
$ws.Range("A28:I28").Select()
